# Insert a new weekly price-report row for "Vega Monumental Concepción - Frutilla"
# at row 229, pushing the existing rows 229-319 down to 230-320.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 229:319 down by one row (native Excel row insert keeps all the
# untouched rows/values/formatting intact and auto-extends the used range).
$ws.Rows("229:229").Insert()

# Populate the newly-inserted row 229 with the new data point.
$ws.Cells.Item(229, 1).Value = 11
$ws.Cells.Item(229, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(229, 3).Value = "Bíobío"
$ws.Cells.Item(229, 4).Value = 44636
$ws.Cells.Item(229, 5).Value = 8
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100101
$ws.Cells.Item(229, 8).Value = "Berries"
$ws.Cells.Item(229, 9).Value = 100112025
$ws.Cells.Item(229, 10).Value = "Frutilla"
$ws.Cells.Item(229, 11).Value = "Sin especificar"
$ws.Cells.Item(229, 12).Value = "Primera"
$ws.Cells.Item(229, 13).Value = 100
$ws.Cells.Item(229, 14).Value = 7000
$ws.Cells.Item(229, 15).Value = 7000
$ws.Cells.Item(229, 16).Value = 7000
$ws.Cells.Item(229, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(229, 18).Value = "Región del Maule"
$ws.Cells.Item(229, 19).Value = 1000
$ws.Cells.Item(229, 20).Value = 7
